# Insert two new price-observation rows right before the existing row 320
# (everything currently at row 320 and below shifts down by two rows), then
# populate the two new rows with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 320 downwards by inserting two blank rows at 320:321.
$ws.Rows("320:321").Insert()

# New row 320: Argentina(o) / Primera
$ws.Range("A320").Value = 3
$ws.Range("B320").Value = "Femacal de La Calera"
$ws.Range("C320").Value = "Coquimbo"
$ws.Range("D320").Value = "2022-08-25"
$ws.Range("E320").Value = 5
$ws.Range("F320").Value = 100112013
$ws.Range("G320").Value = "Alcachofa"
$ws.Range("H320").Value = "Argentina(o)"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 111
$ws.Range("K320").Value = 12500
$ws.Range("L320").Value = 13000
$ws.Range("M320").Value = 12703
$ws.Range("N320").Value = "$/caja 50 unidades"
$ws.Range("O320").Value = "Provincia de Limarí"
$ws.Range("P320").Value = 254
$ws.Range("Q320").Value = 50
$ws.Range("R320").Value = "Hortaliza"

# New row 321: Española / Extra
$ws.Range("A321").Value = 3
$ws.Range("B321").Value = "Femacal de La Calera"
$ws.Range("C321").Value = "Coquimbo"
$ws.Range("D321").Value = "2022-08-25"
$ws.Range("E321").Value = 5
$ws.Range("F321").Value = 100112013
$ws.Range("G321").Value = "Alcachofa"
$ws.Range("H321").Value = "Española"
$ws.Range("I321").Value = "Extra"
$ws.Range("J321").Value = 125
$ws.Range("K321").Value = 14000
$ws.Range("L321").Value = 14500
$ws.Range("M321").Value = 14260
$ws.Range("N321").Value = "$/caja 30 unidades"
$ws.Range("O321").Value = "Provincia de Limarí"
$ws.Range("P321").Value = 475
$ws.Range("Q321").Value = 30
$ws.Range("R321").Value = "Hortaliza"
